# GPLIM-5942 remove Pooled Dev Tube experiment and condition
# Delete columns H:I ("Experiment" and "Conditions") from the sheet,
# shifting all subsequent columns left by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H:I").Delete()
